$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data range we touch (B2:G51) is normally General-formatted text (inlineStr) cells.
# Setting .Value on a numeric-looking string (e.g. "243.44" or "0") would otherwise be
# auto-coerced by Excel into a real number. To keep these as text (matching the source
# file), force NumberFormat to "@" (Text) before assigning, then restore the default
# "Normal" style afterwards so the cells end up with no lingering explicit format -
# same as before the edit, just with new text content.
$dataRange = $ws.Range("B2:G51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "243.44"
$ws.Cells.Item(2, 6).Value = "27-12-2022"
$ws.Cells.Item(2, 7).Value = "0"

# Row 3
$ws.Cells.Item(3, 4).Value = "23.03"
$ws.Cells.Item(3, 6).Value = "27-12-2022"
$ws.Cells.Item(3, 7).Value = "0"

# Row 4
$ws.Cells.Item(4, 4).Value = "5.427"
$ws.Cells.Item(4, 6).Value = "27-12-2022"
$ws.Cells.Item(4, 7).Value = "0"

# Row 5
$ws.Cells.Item(5, 4).Value = "0.05969"
$ws.Cells.Item(5, 6).Value = "27-12-2022"
$ws.Cells.Item(5, 7).Value = "0"

# Row 6
$ws.Cells.Item(6, 4).Value = "3.449"
$ws.Cells.Item(6, 6).Value = "27-12-2022"
$ws.Cells.Item(6, 7).Value = "0"

# Row 7
$ws.Cells.Item(7, 4).Value = "6.527"
$ws.Cells.Item(7, 6).Value = "27-12-2022"
$ws.Cells.Item(7, 7).Value = "0"

# Row 8
$ws.Cells.Item(8, 4).Value = "0.8135"
$ws.Cells.Item(8, 6).Value = "27-12-2022"
$ws.Cells.Item(8, 7).Value = "0"

# Row 9
$ws.Cells.Item(9, 4).Value = "0.9274"
$ws.Cells.Item(9, 6).Value = "27-12-2022"
$ws.Cells.Item(9, 7).Value = "0"

# Row 10
$ws.Cells.Item(10, 4).Value = "0.1427"
$ws.Cells.Item(10, 6).Value = "27-12-2022"
$ws.Cells.Item(10, 7).Value = "0"

# Row 11
$ws.Cells.Item(11, 4).Value = "0.07419"
$ws.Cells.Item(11, 6).Value = "27-12-2022"
$ws.Cells.Item(11, 7).Value = "0"

# Row 12
$ws.Cells.Item(12, 4).Value = "0.03299"
$ws.Cells.Item(12, 6).Value = "27-12-2022"
$ws.Cells.Item(12, 7).Value = "0"

# Row 13
$ws.Cells.Item(13, 4).Value = "0.03094"
$ws.Cells.Item(13, 6).Value = "27-12-2022"
$ws.Cells.Item(13, 7).Value = "0"

# Row 14
$ws.Cells.Item(14, 4).Value = "0.09365"
$ws.Cells.Item(14, 6).Value = "27-12-2022"
$ws.Cells.Item(14, 7).Value = "0"

# Row 15
$ws.Cells.Item(15, 4).Value = "3.859"
$ws.Cells.Item(15, 6).Value = "27-12-2022"
$ws.Cells.Item(15, 7).Value = "0"

# Row 16
$ws.Cells.Item(16, 6).Value = "27-12-2022"
$ws.Cells.Item(16, 7).Value = "0"

# Row 17
$ws.Cells.Item(17, 4).Value = "0.04697"
$ws.Cells.Item(17, 6).Value = "27-12-2022"
$ws.Cells.Item(17, 7).Value = "0"

# Row 18
$ws.Cells.Item(18, 4).Value = "0.0005896"
$ws.Cells.Item(18, 5).Value = "17OneONEWorstin24h"
$ws.Cells.Item(18, 6).Value = "27-12-2022"
$ws.Cells.Item(18, 7).Value = "0"

# Row 19
$ws.Cells.Item(19, 4).Value = "0.005859"
$ws.Cells.Item(19, 6).Value = "27-12-2022"
$ws.Cells.Item(19, 7).Value = "0"

# Row 20
$ws.Cells.Item(20, 4).Value = "0.001262"
$ws.Cells.Item(20, 5).Value = "19BitKanKAN"
$ws.Cells.Item(20, 6).Value = "27-12-2022"
$ws.Cells.Item(20, 7).Value = "0"

# Row 21
$ws.Cells.Item(21, 4).Value = "0.004916"
$ws.Cells.Item(21, 6).Value = "27-12-2022"
$ws.Cells.Item(21, 7).Value = "0"

# Row 22
$ws.Cells.Item(22, 4).Value = "0.00006805"
$ws.Cells.Item(22, 6).Value = "27-12-2022"
$ws.Cells.Item(22, 7).Value = "0"

# Row 23
$ws.Cells.Item(23, 4).Value = "3.563"
$ws.Cells.Item(23, 6).Value = "27-12-2022"
$ws.Cells.Item(23, 7).Value = "0"

# Row 24
$ws.Cells.Item(24, 4).Value = "2.154"
$ws.Cells.Item(24, 6).Value = "27-12-2022"
$ws.Cells.Item(24, 7).Value = "0"

# Row 25
$ws.Cells.Item(25, 4).Value = "0.3233"
$ws.Cells.Item(25, 6).Value = "27-12-2022"
$ws.Cells.Item(25, 7).Value = "0"

# Row 26
$ws.Cells.Item(26, 4).Value = "0.1333"
$ws.Cells.Item(26, 6).Value = "27-12-2022"
$ws.Cells.Item(26, 7).Value = "0"

# Row 27
$ws.Cells.Item(27, 4).Value = "0.0002304"
$ws.Cells.Item(27, 6).Value = "27-12-2022"
$ws.Cells.Item(27, 7).Value = "0"

# Row 28
$ws.Cells.Item(28, 6).Value = "27-12-2022"
$ws.Cells.Item(28, 7).Value = "0"

# Row 29
$ws.Cells.Item(29, 6).Value = "27-12-2022"
$ws.Cells.Item(29, 7).Value = "0"

# Row 30
$ws.Cells.Item(30, 6).Value = "27-12-2022"
$ws.Cells.Item(30, 7).Value = "0"

# Row 31
$ws.Cells.Item(31, 6).Value = "27-12-2022"
$ws.Cells.Item(31, 7).Value = "0"

# Row 32
$ws.Cells.Item(32, 6).Value = "27-12-2022"
$ws.Cells.Item(32, 7).Value = "0"

# Row 33
$ws.Cells.Item(33, 6).Value = "27-12-2022"
$ws.Cells.Item(33, 7).Value = "0"

# Row 34
$ws.Cells.Item(34, 6).Value = "27-12-2022"
$ws.Cells.Item(34, 7).Value = "0"

# Row 35
$ws.Cells.Item(35, 6).Value = "27-12-2022"
$ws.Cells.Item(35, 7).Value = "0"

# Row 36
$ws.Cells.Item(36, 6).Value = "27-12-2022"
$ws.Cells.Item(36, 7).Value = "0"

# Row 37
$ws.Cells.Item(37, 6).Value = "27-12-2022"
$ws.Cells.Item(37, 7).Value = "0"

# Row 38
$ws.Cells.Item(38, 6).Value = "27-12-2022"
$ws.Cells.Item(38, 7).Value = "0"

# Row 39
$ws.Cells.Item(39, 6).Value = "27-12-2022"
$ws.Cells.Item(39, 7).Value = "0"

# Row 40
$ws.Cells.Item(40, 4).Value = "0.03968"
$ws.Cells.Item(40, 6).Value = "27-12-2022"
$ws.Cells.Item(40, 7).Value = "0"

# Row 41
$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(41, 4).Value = "0.006465"
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"
$ws.Cells.Item(41, 6).Value = "27-12-2022"
$ws.Cells.Item(41, 7).Value = "0"

# Row 42
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(42, 4).Value = "0.1077"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"
$ws.Cells.Item(42, 6).Value = "27-12-2022"
$ws.Cells.Item(42, 7).Value = "0"

# Row 43
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(43, 4).Value = "0.003003"
$ws.Cells.Item(43, 5).Value = "42CEJICEJI"
$ws.Cells.Item(43, 6).Value = "27-12-2022"
$ws.Cells.Item(43, 7).Value = "0"

# Row 44
$ws.Cells.Item(44, 4).Value = "0.009008"
$ws.Cells.Item(44, 5).Value = "43LocalTradersLCTBestin24h"
$ws.Cells.Item(44, 6).Value = "27-12-2022"
$ws.Cells.Item(44, 7).Value = "0"

# Row 45
$ws.Cells.Item(45, 4).Value = "0.00005239"
$ws.Cells.Item(45, 6).Value = "27-12-2022"
$ws.Cells.Item(45, 7).Value = "0"

# Row 46
$ws.Cells.Item(46, 4).Value = "0.00000000751"
$ws.Cells.Item(46, 6).Value = "27-12-2022"
$ws.Cells.Item(46, 7).Value = "0"

# Row 47
$ws.Cells.Item(47, 4).Value = "0.6707"
$ws.Cells.Item(47, 6).Value = "27-12-2022"
$ws.Cells.Item(47, 7).Value = "0"

# Row 48
$ws.Cells.Item(48, 4).Value = "0.002394"
$ws.Cells.Item(48, 6).Value = "27-12-2022"
$ws.Cells.Item(48, 7).Value = "0"

# Row 49
$ws.Cells.Item(49, 4).Value = "0.00002102"
$ws.Cells.Item(49, 6).Value = "27-12-2022"
$ws.Cells.Item(49, 7).Value = "0"

# Row 50
$ws.Cells.Item(50, 4).Value = "0.0002002"
$ws.Cells.Item(50, 6).Value = "27-12-2022"
$ws.Cells.Item(50, 7).Value = "0"

# Row 51
$ws.Cells.Item(51, 6).Value = "27-12-2022"
$ws.Cells.Item(51, 7).Value = "0"

# Restore the default (General) style on the whole data range so only cell content
# changed - no stray "@" text-format styling is left behind.
$dataRange.Style = "Normal"